$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 33
$ws.Cells.Item(2, 3).Value = 'Arts for Lawrence'
$ws.Cells.Item(2, 4).Value = 4.9
$ws.Cells.Item(2, 5).Value = 19
$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 3).Value = 'Box Office at The Center for Performing Arts'
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(4, 1).Value = 8
$ws.Cells.Item(4, 3).Value = 'Carmel Gazebo'
$ws.Cells.Item(4, 4).Value = 4.6
$ws.Cells.Item(4, 5).Value = 54
$ws.Cells.Item(5, 1).Value = 36
$ws.Cells.Item(5, 3).Value = 'Celebration Plaza Amphitheater'
$ws.Cells.Item(5, 4).Value = 4.6
$ws.Cells.Item(5, 5).Value = 26
$ws.Cells.Item(6, 1).Value = 51
$ws.Cells.Item(6, 3).Value = 'Circle City Tickets'
$ws.Cells.Item(6, 4).Value = 4.5
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(7, 1).Value = 27
$ws.Cells.Item(7, 3).Value = 'Clowes Memorial Hall'
$ws.Cells.Item(7, 4).Value = 4.7
$ws.Cells.Item(7, 5).Value = 797
$ws.Cells.Item(8, 1).Value = 58
$ws.Cells.Item(8, 3).Value = 'Conner Prairie'
$ws.Cells.Item(8, 4).Value = 4.6
$ws.Cells.Item(8, 5).Value = 1570
$ws.Cells.Item(9, 1).Value = 18
$ws.Cells.Item(9, 3).Value = 'Convention center'
$ws.Cells.Item(9, 4).Value = 4.6
$ws.Cells.Item(9, 5).Value = 473
$ws.Cells.Item(10, 1).Value = 35
$ws.Cells.Item(10, 3).Value = 'Cool Creek Park Nature Center'
$ws.Cells.Item(10, 4).Value = 4.7
$ws.Cells.Item(10, 5).Value = 462
$ws.Cells.Item(11, 1).Value = 12
$ws.Cells.Item(11, 3).Value = 'Cool Creek Park Open Theatre/Stadium'
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(12, 1).Value = 49
$ws.Cells.Item(12, 3).Value = 'Dr. James A. Dillon Park'
$ws.Cells.Item(12, 4).Value = 4.6
$ws.Cells.Item(12, 5).Value = 461
$ws.Cells.Item(13, 1).Value = 29
$ws.Cells.Item(13, 3).Value = 'Egyptian Room at Old National Centre'
$ws.Cells.Item(13, 4).Value = 4.5
$ws.Cells.Item(13, 5).Value = 41
$ws.Cells.Item(14, 1).Value = 10
$ws.Cells.Item(14, 3).Value = 'Eidson-Duckwall Recital Hall'
$ws.Cells.Item(14, 4).Value = 4.6
$ws.Cells.Item(14, 5).Value = 19
$ws.Cells.Item(15, 1).Value = 25
$ws.Cells.Item(15, 3).Value = 'Elvis Presley Final Concert Plaque'
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(16, 1).Value = 5
$ws.Cells.Item(16, 3).Value = 'Emerson Theater'
$ws.Cells.Item(16, 4).Value = 4
$ws.Cells.Item(16, 5).Value = 491
$ws.Cells.Item(17, 1).Value = 21
$ws.Cells.Item(17, 3).Value = 'Federal Hill Commons'
$ws.Cells.Item(17, 4).Value = 4.6
$ws.Cells.Item(17, 5).Value = 455
$ws.Cells.Item(18, 1).Value = 55
$ws.Cells.Item(18, 3).Value = 'Fountain Square Theatre'
$ws.Cells.Item(18, 4).Value = 4.6
$ws.Cells.Item(18, 5).Value = 167
$ws.Cells.Item(19, 1).Value = 19
$ws.Cells.Item(19, 3).Value = 'Free Stage'
$ws.Cells.Item(19, 4).Value = 4.6
$ws.Cells.Item(19, 5).Value = 361
$ws.Cells.Item(20, 1).Value = 41
$ws.Cells.Item(20, 3).Value = 'Gainbridge Fieldhouse'
$ws.Cells.Item(20, 4).Value = 4.7
$ws.Cells.Item(20, 5).Value = 9144
$ws.Cells.Item(21, 1).Value = 15
$ws.Cells.Item(21, 3).Value = 'HI-FI Indy & HI-FI Annex'
$ws.Cells.Item(21, 4).Value = 4.6
$ws.Cells.Item(21, 5).Value = 1113
$ws.Cells.Item(22, 1).Value = 11
$ws.Cells.Item(22, 3).Value = 'Hilbert Circle Theatre'
$ws.Cells.Item(22, 4).Value = 4.8
$ws.Cells.Item(22, 5).Value = 978
$ws.Cells.Item(23, 1).Value = 14
$ws.Cells.Item(23, 3).Value = 'Hoosier Dome'
$ws.Cells.Item(23, 4).Value = 4.4
$ws.Cells.Item(23, 5).Value = 226
$ws.Cells.Item(24, 1).Value = 47
$ws.Cells.Item(24, 3).Value = 'IMMI Conference Center'
$ws.Cells.Item(24, 4).Value = 4.6
$ws.Cells.Item(24, 5).Value = 23
$ws.Cells.Item(25, 1).Value = 37
$ws.Cells.Item(25, 3).Value = 'Indiana Historical Society'
$ws.Cells.Item(25, 4).Value = 4.7
$ws.Cells.Item(25, 5).Value = 730
$ws.Cells.Item(26, 1).Value = 48
$ws.Cells.Item(26, 3).Value = 'Indiana State Fairgrounds & Event Center'
$ws.Cells.Item(26, 4).Value = 4.4
$ws.Cells.Item(26, 5).Value = 1344
$ws.Cells.Item(27, 1).Value = 46
$ws.Cells.Item(27, 3).Value = 'Indianapolis Chamber Orchestra'
$ws.Cells.Item(27, 4).Value = 5
$ws.Cells.Item(27, 5).Value = 5
$ws.Cells.Item(28, 1).Value = 52
$ws.Cells.Item(28, 3).Value = 'Indianapolis Motor Speedway'
$ws.Cells.Item(28, 4).Value = 4.8
$ws.Cells.Item(28, 5).Value = 11013
$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 3).Value = 'Indianapolis Symphony Orchestra'
$ws.Cells.Item(29, 4).Value = 4.1
$ws.Cells.Item(29, 5).Value = 8
$ws.Cells.Item(30, 1).Value = 57
$ws.Cells.Item(30, 3).Value = 'Indianapolis Zoo'
$ws.Cells.Item(30, 4).Value = 4.5
$ws.Cells.Item(30, 5).Value = 14982
$ws.Cells.Item(31, 1).Value = 59
$ws.Cells.Item(31, 3).Value = 'IndyFringe Theatre'
$ws.Cells.Item(31, 4).Value = 4.6
$ws.Cells.Item(31, 5).Value = 170
$ws.Cells.Item(32, 1).Value = 28
$ws.Cells.Item(32, 3).Value = 'Irving Theater'
$ws.Cells.Item(32, 4).Value = 4.3
$ws.Cells.Item(32, 5).Value = 382
$ws.Cells.Item(33, 1).Value = 42
$ws.Cells.Item(33, 3).Value = 'Lincoln Park'
$ws.Cells.Item(33, 4).Value = 4.7
$ws.Cells.Item(33, 5).Value = 21
$ws.Cells.Item(34, 1).Value = 34
$ws.Cells.Item(34, 3).Value = 'Live Nation'
$ws.Cells.Item(34, 4).Value = 3
$ws.Cells.Item(34, 5).Value = 25
$ws.Cells.Item(35, 1).Value = 54
$ws.Cells.Item(35, 3).Value = 'Lucas Oil Stadium'
$ws.Cells.Item(35, 4).Value = 4.7
$ws.Cells.Item(35, 5).Value = 11722
$ws.Cells.Item(36, 1).Value = 20
$ws.Cells.Item(36, 3).Value = 'MOKB Presents'
$ws.Cells.Item(36, 4).Value = 3.6
$ws.Cells.Item(36, 5).Value = 5
$ws.Cells.Item(37, 1).Value = 53
$ws.Cells.Item(37, 3).Value = 'McGowan Hall'
$ws.Cells.Item(37, 4).Value = 4.6
$ws.Cells.Item(37, 5).Value = 116
$ws.Cells.Item(38, 1).Value = 40
$ws.Cells.Item(38, 3).Value = 'Meadowood Park'
$ws.Cells.Item(38, 4).Value = 4.6
$ws.Cells.Item(38, 5).Value = 702
$ws.Cells.Item(39, 1).Value = 31
$ws.Cells.Item(39, 3).Value = 'Military Park'
$ws.Cells.Item(39, 4).Value = 4.6
$ws.Cells.Item(39, 5).Value = 1512
$ws.Cells.Item(40, 1).Value = 44
$ws.Cells.Item(40, 3).Value = 'Mimi’s Event Lounge/Moe’s Cafe'
$ws.Cells.Item(40, 4).Value = 3.9
$ws.Cells.Item(40, 5).Value = 52
$ws.Cells.Item(41, 1).Value = 38
$ws.Cells.Item(41, 3).Value = 'Murat Egyptian Room'
$ws.Cells.Item(41, 4).Value = 4.6
$ws.Cells.Item(41, 5).Value = 65
$ws.Cells.Item(42, 1).Value = 22
$ws.Cells.Item(42, 3).Value = 'Murat Theatre'
$ws.Cells.Item(42, 4).Value = 4.5
$ws.Cells.Item(42, 5).Value = 1292
$ws.Cells.Item(43, 1).Value = 1
$ws.Cells.Item(43, 3).Value = 'Nickel Plate District Amphitheater'
$ws.Cells.Item(43, 4).Value = 4.6
$ws.Cells.Item(43, 5).Value = 430
$ws.Cells.Item(44, 1).Value = 3
$ws.Cells.Item(44, 3).Value = 'Old National Centre'
$ws.Cells.Item(44, 4).Value = 4.4
$ws.Cells.Item(44, 5).Value = 2421
$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 3).Value = 'Philharmonic Orchestra'
$ws.Cells.Item(45, 4).Value = 5
$ws.Cells.Item(45, 5).Value = 1
$ws.Cells.Item(46, 1).Value = 30
$ws.Cells.Item(46, 3).Value = 'Pike Performing Arts Center'
$ws.Cells.Item(46, 4).Value = 4.6
$ws.Cells.Item(46, 5).Value = 180
$ws.Cells.Item(47, 1).Value = 13
$ws.Cells.Item(47, 3).Value = 'Ruoff Home Mortgage center'
$ws.Cells.Item(47, 4).Value = 4.4
$ws.Cells.Item(47, 5).Value = 16
$ws.Cells.Item(48, 1).Value = 32
$ws.Cells.Item(48, 3).Value = 'Sidewalk Concert'
$ws.Cells.Item(48, 4).Value = 0
$ws.Cells.Item(48, 5).Value = 0
$ws.Cells.Item(49, 1).Value = 2
$ws.Cells.Item(49, 3).Value = 'TCU Amphitheater at White River State Park'
$ws.Cells.Item(49, 4).Value = 4.6
$ws.Cells.Item(49, 5).Value = 1188
$ws.Cells.Item(50, 1).Value = 50
$ws.Cells.Item(50, 3).Value = 'The Cabaret'
$ws.Cells.Item(50, 4).Value = 4.8
$ws.Cells.Item(50, 5).Value = 90
$ws.Cells.Item(51, 1).Value = 9
$ws.Cells.Item(51, 3).Value = 'The Center for the Performing Arts'
$ws.Cells.Item(51, 4).Value = 4.7
$ws.Cells.Item(51, 5).Value = 1193
$ws.Cells.Item(52, 1).Value = 16
$ws.Cells.Item(52, 3).Value = 'The Crane Bay Event Center'
$ws.Cells.Item(52, 4).Value = 4.6
$ws.Cells.Item(52, 5).Value = 248
$ws.Cells.Item(53, 1).Value = 17
$ws.Cells.Item(53, 3).Value = 'The Palladium at the Center for the Performing Arts'
$ws.Cells.Item(53, 4).Value = 4.9
$ws.Cells.Item(53, 5).Value = 197
$ws.Cells.Item(54, 1).Value = 23
$ws.Cells.Item(54, 3).Value = 'The Pavilion at Pan Am'
$ws.Cells.Item(54, 4).Value = 4.4
$ws.Cells.Item(54, 5).Value = 460
$ws.Cells.Item(55, 1).Value = 7
$ws.Cells.Item(55, 3).Value = 'The Vogue'
$ws.Cells.Item(55, 4).Value = 4.5
$ws.Cells.Item(55, 5).Value = 1396
$ws.Cells.Item(56, 1).Value = 56
$ws.Cells.Item(56, 3).Value = 'Theater at the Fort'
$ws.Cells.Item(56, 4).Value = 4.5
$ws.Cells.Item(56, 5).Value = 75
$ws.Cells.Item(57, 1).Value = 24
$ws.Cells.Item(57, 3).Value = 'Ticketmaster'
$ws.Cells.Item(57, 4).Value = 2
$ws.Cells.Item(57, 5).Value = 16
$ws.Cells.Item(58, 1).Value = 4
$ws.Cells.Item(58, 3).Value = 'Warren Performing Arts Center'
$ws.Cells.Item(58, 4).Value = 4.5
$ws.Cells.Item(58, 5).Value = 374
$ws.Cells.Item(59, 1).Value = 39
$ws.Cells.Item(59, 3).Value = 'Washington Township Park'
$ws.Cells.Item(59, 4).Value = 4.6
$ws.Cells.Item(59, 5).Value = 1018
$ws.Cells.Item(60, 1).Value = 45
$ws.Cells.Item(60, 3).Value = 'White River State Park'
$ws.Cells.Item(60, 4).Value = 4.7
$ws.Cells.Item(60, 5).Value = 3838
$ws.Cells.Item(61, 1).Value = 26
$ws.Cells.Item(61, 3).Value = 'Zionsville Concert Band'
$ws.Cells.Item(61, 4).Value = 0
$ws.Cells.Item(61, 5).Value = 0
